$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Replace the paragraph's text (drop the trailing run with just a space,
# and update the bookmark-style placeholder text) while keeping the
# first run's character formatting.
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "**ID__AFFARS_AFMC_PGI_5315_3A__ID**"

# Re-fetch the (now single-run) paragraph and update its paragraph
# formatting: add a paragraph border (5pt space on every side, no line)
# and widen the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p = $d.Paragraphs(1)
$p.Format.LeftIndent = 11.25
$p.Range.Borders.DistanceFromTop = 5
$p.Range.Borders.DistanceFromLeft = 5
$p.Range.Borders.DistanceFromBottom = 5
$p.Range.Borders.DistanceFromRight = 5

Write-Output "done"
